# Auto-generated edit script applying numeric corrections per sheet/row/column
# as described in the commit diff (scheduled runner price-refresh update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2516.6667
$ws.Range("I62").Value = 2516.6667
$ws.Range("K62").Value = 2516.6667
$ws.Range("M62").Value = -1892.6667

# Row 65
$ws.Range("H65").Value = 2516.6667
$ws.Range("I65").Value = 2516.6667
$ws.Range("K65").Value = 12583.3335
$ws.Range("M65").Value = -9463.3335

# Row 132
$ws.Range("H132").Value = 6950058
$ws.Range("I132").Value = 9525566
$ws.Range("K132").Value = 28576698
$ws.Range("M132").Value = -28574168

# Row 137
$ws.Range("H137").Value = 1163.8572
$ws.Range("I137").Value = 1144
$ws.Range("J137").Value = 1248.25
$ws.Range("K137").Value = 3432
$ws.Range("L137").Value = 3744.75
$ws.Range("M137").Value = -882
$ws.Range("N137").Value = -8844.75

# Row 138
$ws.Range("H138").Value = 580864.3
$ws.Range("I138").Value = 634.0263
$ws.Range("J138").Value = 1176776.5
$ws.Range("K138").Value = 1902.0789
$ws.Range("L138").Value = 3530329.5
$ws.Range("M138").Value = 3237.9211
$ws.Range("N138").Value = -3540609.5

# Row 141
$ws.Range("H141").Value = 1095
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 5493.636
$ws.Range("I2").Value = 993.17645
$ws.Range("J2").Value = 20795.2
$ws.Range("K2").Value = 993.17645
$ws.Range("L2").Value = 20795.2
$ws.Range("M2").Value = -880.17645
$ws.Range("N2").Value = -21021.2

# Row 61
$ws.Range("H61").Value = 1235.5
$ws.Range("I61").Value = 1099.7693
$ws.Range("K61").Value = 1099.7693
$ws.Range("M61").Value = -887.7692999999999

# Row 74
$ws.Range("H74").Value = 760.5455
$ws.Range("I74").Value = 760.5455
$ws.Range("K74").Value = 760.5455
$ws.Range("M74").Value = 113.4545000000001

# Row 77
$ws.Range("H77").Value = 760.5455
$ws.Range("I77").Value = 760.5455
$ws.Range("K77").Value = 3802.7275
$ws.Range("M77").Value = 565.2725

# Row 116
$ws.Range("H116").Value = 5493.636
$ws.Range("I116").Value = 993.17645
$ws.Range("J116").Value = 20795.2
$ws.Range("K116").Value = 993.17645
$ws.Range("L116").Value = 20795.2
$ws.Range("M116").Value = 1300.82355
$ws.Range("N116").Value = -25383.2

# Row 136
$ws.Range("H136").Value = 1235.5
$ws.Range("I136").Value = 1099.7693
$ws.Range("K136").Value = 3299.3079
$ws.Range("M136").Value = -749.3078999999998


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 5493.636
$ws.Range("I3").Value = 993.17645
$ws.Range("J3").Value = 20795.2
$ws.Range("K3").Value = 993.17645
$ws.Range("L3").Value = 20795.2
$ws.Range("M3").Value = -879.17645
$ws.Range("N3").Value = -21023.2

# Row 98
$ws.Range("H98").Value = 75000
$ws.Range("J98").Value = 75000
$ws.Range("L98").Value = 75000
$ws.Range("N98").Value = -80990

# Row 134
$ws.Range("H134").Value = 4778.607
$ws.Range("I134").Value = 1251.9259
$ws.Range("K134").Value = 3755.7777
$ws.Range("M134").Value = -1220.7777


$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 583
$ws.Range("J22").Value = 666
$ws.Range("L22").Value = 666
$ws.Range("N22").Value = -1366

# Row 31
$ws.Range("H31").Value = 2070.795
$ws.Range("I31").Value = 1004.7857
$ws.Range("J31").Value = 2667.76
$ws.Range("K31").Value = 1004.7857
$ws.Range("L31").Value = 2667.76
$ws.Range("M31").Value = -709.7857
$ws.Range("N31").Value = -3257.76

# Row 34
$ws.Range("H34").Value = 2070.795
$ws.Range("I34").Value = 1004.7857
$ws.Range("J34").Value = 2667.76
$ws.Range("K34").Value = 1004.7857
$ws.Range("L34").Value = 2667.76
$ws.Range("M34").Value = -802.7857
$ws.Range("N34").Value = -3071.76

# Row 132
$ws.Range("H132").Value = 1887.4783
$ws.Range("I132").Value = 1103.9286
$ws.Range("J132").Value = 3106.3333
$ws.Range("K132").Value = 3311.7858
$ws.Range("L132").Value = 9318.999899999999
$ws.Range("M132").Value = -781.7857999999997
$ws.Range("N132").Value = -14378.9999


$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("K16").Value = 300
$ws.Range("M16").Value = -127

# Row 116
$ws.Range("H116").Value = 3619.2
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# Row 131
$ws.Range("H131").Value = 14494011
$ws.Range("I131").Value = 166667090
$ws.Range("J131").Value = 1337.4921
$ws.Range("K131").Value = 500001270
$ws.Range("L131").Value = 4012.4763
$ws.Range("M131").Value = -499996230
$ws.Range("N131").Value = -14092.4763


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2220.8
$ws.Range("I7").Value = 2199.6667
$ws.Range("J7").Value = 2252.5
$ws.Range("K7").Value = 2199.6667
$ws.Range("L7").Value = 2252.5
$ws.Range("M7").Value = -2087.6667
$ws.Range("N7").Value = -2476.5

# Row 40
$ws.Range("H40").Value = 3197.818
$ws.Range("I40").Value = 2834.5
$ws.Range("J40").Value = 4166.6665
$ws.Range("K40").Value = 2834.5
$ws.Range("L40").Value = 4166.6665
$ws.Range("M40").Value = -2698.5
$ws.Range("N40").Value = -4438.6665

# Row 46
$ws.Range("H46").Value = 1701.6428
$ws.Range("J46").Value = 1861.6
$ws.Range("L46").Value = 1861.6
$ws.Range("N46").Value = -2237.6

# Row 61
$ws.Range("H61").Value = 1448.1428
$ws.Range("I61").Value = 1489.5
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 1489.5
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -1287.5
$ws.Range("N61").Value = -1604

# Row 68
$ws.Range("H68").Value = 1677.5883
$ws.Range("I68").Value = 1608.7693
$ws.Range("K68").Value = 1608.7693
$ws.Range("M68").Value = -859.7692999999999

# Row 71
$ws.Range("H71").Value = 1677.5883
$ws.Range("I71").Value = 1608.7693
$ws.Range("K71").Value = 8043.8465
$ws.Range("M71").Value = -4299.8465

# Row 113
$ws.Range("H113").Value = 1448.1428
$ws.Range("I113").Value = 1489.5
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1489.5
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 680.5
$ws.Range("N113").Value = -5540

# Row 126
$ws.Range("H126").Value = 2220.8
$ws.Range("I126").Value = 2199.6667
$ws.Range("J126").Value = 2252.5
$ws.Range("K126").Value = 6599.000100000001
$ws.Range("L126").Value = 6757.5
$ws.Range("M126").Value = -4129.000100000001
$ws.Range("N126").Value = -11697.5

# Row 132
$ws.Range("H132").Value = 64552.062
$ws.Range("I132").Value = 1271.4286
$ws.Range("J132").Value = 113770.336
$ws.Range("K132").Value = 3814.2858
$ws.Range("L132").Value = 341311.008
$ws.Range("M132").Value = -1284.2858
$ws.Range("N132").Value = -346371.008


$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 530.1429
$ws.Range("I107").Value = 548.75
$ws.Range("K107").Value = 1646.25
$ws.Range("M107").Value = 273.75

# Row 126
$ws.Range("H126").Value = 123457450
$ws.Range("I126").Value = 158730740
$ws.Range("K126").Value = 476192220
$ws.Range("M126").Value = -476189750

# Row 132
$ws.Range("H132").Value = 2925.375
$ws.Range("I132").Value = 2593.9
$ws.Range("K132").Value = 7781.700000000001
$ws.Range("M132").Value = -5251.700000000001

